$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "Save" header in H1, matching the style/formatting of the other
# header cells (copy G1's format into H1, then overwrite the value/text)
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Add the new "Save" column values for the data rows (plain numeric cells,
# same as F2/F3/G2/G3 which carry no explicit style)
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
